$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.898.61"
$ws.Range("E2").Value = "  -1.65%  "
$ws.Range("D3").Value = "1.887.83"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.78%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7351"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.79"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3109"
$ws.Range("E8").Value = "  -2.40%  "
$ws.Range("E9").Value = "  -5.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06901"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7730"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07942"
$ws.Range("D13").Value = "1.873.27"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.224"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.30"
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.21"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "29.919.80"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.760"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  -5.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007755"
$ws.Range("E20").Value = "  -1.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "2.126.32"
$ws.Range("E22").Value = "  -2.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.918"
$ws.Range("E24").Value = "  +3.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.303"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.33"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.86"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E28").Value = "  -4.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.019"
$ws.Range("E29").Value = "  -10.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.352"
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.532"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.303"
$ws.Range("E32").Value = "  -1.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.081"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05096"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.277"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7365"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01919"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.296"
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.31"
$ws.Range("E41").Value = "  -5.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4462"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.931"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8369"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.640"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.82"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.775"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "2.051.75"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.92"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "942.05"
$ws.Range("E51").Value = "  -3.14%  "
